# Update PLC data 2025-10-13 13:38:14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 252
$ws.Range("C3").Value = 152796
$ws.Range("C4").Value = 144358
$ws.Range("C7").Value = 5.52
$ws.Range("C8").Value = 63.65
